$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) values are written as text, preserving exact formatting
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "66.779.78"
$ws.Range("E2").Value = "  -0.12%  "
$ws.Range("D3").Value = "3.447.67"
$ws.Range("E3").Value = "  -0.72%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "585.98"
$ws.Range("E5").Value = "  +0.13%  "
$ws.Range("D6").Value = "179.29"
$ws.Range("E6").Value = "  +1.69%  "
$ws.Range("D7").Value = "0.632"
$ws.Range("E7").Value = "  +5.53%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").Value = "3.445.48"
$ws.Range("E9").Value = "  -0.76%  "
$ws.Range("E10").Value = "  +0.69%  "
$ws.Range("D11").Value = "6.97"
$ws.Range("E11").Value = "  +0.98%  "
$ws.Range("D12").Value = "0.419"
$ws.Range("E12").Value = "  -0.78%  "
$ws.Range("D13").Value = "4.047.38"
$ws.Range("E13").Value = "  -0.61%  "
$ws.Range("E14").Value = "  +1.41%  "
$ws.Range("D15").Value = "29.93"
$ws.Range("E15").Value = "  -1.53%  "
$ws.Range("D16").Value = "66.723.63"
$ws.Range("E16").Value = "  -0.20%  "
$ws.Range("D17").Value = "0.0000175"
$ws.Range("E17").Value = "  +1.08%  "
$ws.Range("D18").Value = "3.472.17"
$ws.Range("E18").Value = "  +0.41%  "
$ws.Range("D19").Value = "5.97"
$ws.Range("E19").Value = "  -1.14%  "
$ws.Range("D20").Value = "13.89"
$ws.Range("E20").Value = "  -0.05%  "
$ws.Range("D21").Value = "373.14"
$ws.Range("E21").Value = "  -2.38%  "
$ws.Range("D22").Value = "7.69"
$ws.Range("E22").Value = "  -2.02%  "
$ws.Range("D23").Value = "73.55"
$ws.Range("E23").Value = "  +1.20%  "
$ws.Range("D24").Value = "0.0000131"
$ws.Range("E24").Value = "  +8.61%  "
$ws.Range("E25").Value = "  +0.15%  "
$ws.Range("E26").Value = "  -1.76%  "
$ws.Range("E27").Value = "  +1.19%  "
$ws.Range("E28").Value = "  +2.09%  "
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  +0.02%  "
$ws.Range("E30").Value = "  +0.51%  "
$ws.Range("E31").Value = "  -0.18%  "
$ws.Range("D32").Value = "23.74"
$ws.Range("E32").Value = "  -2.33%  "
$ws.Range("D33").Value = "0.999"
$ws.Range("E33").Value = "  -0.07%  "
$ws.Range("B34").Value = "Fetch.AI"
$ws.Range("C34").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D34").Value = "1.29"
$ws.Range("E34").Value = "  -3.46%  "
$ws.Range("B35").Value = "Aptos"
$ws.Range("C35").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D35").Value = "7.12"
$ws.Range("E35").Value = "  -0.88%  "
$ws.Range("D36").Value = "1.59"
$ws.Range("E36").Value = "  +0.18%  "
$ws.Range("D37").Value = "163.00"
$ws.Range("E37").Value = "  +1.46%  "
$ws.Range("D38").Value = "0.885"
$ws.Range("E38").Value = "  -0.97%  "
$ws.Range("D39").Value = "27.88"
$ws.Range("E39").Value = "  -6.04%  "
$ws.Range("E40").Value = "  +1.25%  "
$ws.Range("D41").Value = "2.64"
$ws.Range("E41").Value = "  +2.27%  "
$ws.Range("D42").Value = "4.52"
$ws.Range("E42").Value = "  +0.05%  "
$ws.Range("D43").Value = "2.758.20"
$ws.Range("E43").Value = "  +2.79%  "
$ws.Range("D44").Value = "6.43"
$ws.Range("E44").Value = "  -0.32%  "
$ws.Range("D45").Value = "0.0700"
$ws.Range("E45").Value = "  +0.02%  "
$ws.Range("D46").Value = "25.72"
$ws.Range("E46").Value = "  +4.28%  "
$ws.Range("D47").Value = "40.10"
$ws.Range("E47").Value = "  -1.30%  "
$ws.Range("D48").Value = "337.93"
$ws.Range("E48").Value = "  +6.85%  "
$ws.Range("D49").Value = "0.0289"
$ws.Range("E49").Value = "  -1.81%  "
$ws.Range("E50").Value = "  +2.93%  "
$ws.Range("D51").Value = "32.03"
$ws.Range("E51").Value = "  +3.80%  "

# Restore default style on Price column (removes temporary text-format style marker)
$priceRange.Style = "Normal"
